$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("L2").Value = "2017-02-15 06:06:40"
$wsZhCn.Range("R2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZhCn.Columns.Item(18).ColumnWidth = 13.7470531463623

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("L2").Value = "2017-02-15 06:07:07"
$wsDeDe.Range("R2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDeDe.Columns.Item(18).ColumnWidth = 13.7470531463623
